# plan.xlsx edit: reschedule posts in rows 16-38 (column F, "Scheduled Time")
# from 2024-04-14 21:34 to 2024-04-15 22:04, and scroll the sheet view down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Scheduled Time" column (F) for rows 16 through 38 -------------
# Old value (all rows): 45396.898611111108  (2024-04-14 21:34:00)
# New value (all rows): 45397.919444444444  (2024-04-15 22:04:00)
$newScheduledTime = 45397.919444444444

for ($row = 16; $row -le 38; $row++) {
    $ws.Cells.Item($row, 6).Value = $newScheduledTime
}

# --- Reposition the sheet view / window scroll ---------------------------------
# The author scrolled the sheet so row 15 is the top visible row (was row 9),
# and resized/repositioned the Excel window.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1

$win.Left = 1886
$win.Top = 1886
$win.Width = 24685
$win.Height = 13097
